$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D) updates ---
# Most target strings (e.g. "36.005.23") contain two dots and are not
# recognized by Excel as numbers, so a plain .Value assignment keeps them
# as text, exactly like the source workbook.
$ws.Range("D2").Value = "36.005.23"
$ws.Range("D3").Value = "1.958.11"
$ws.Range("D16").Value = "2.243.33"
$ws.Range("D18").Value = "1.954.80"
$ws.Range("D19").Value = "35.857.13"
$ws.Range("D49").Value = "1.332.48"
$ws.Range("D51").Value = "2.138.01"

# Some target strings (e.g. "240.85") DO look like plain numbers, and a
# plain .Value assignment would make Excel silently convert them to a
# numeric cell (losing the original text formatting / introducing float
# rounding). To keep them as text without touching any cell style, write a
# formula that evaluates to the literal text, then convert that formula to
# its value in place via Copy + PasteSpecial (values only).
$ws.Range("D5").Formula = "=""240.85"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("D6").Formula = "=""0.622"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("D7").Formula = "=""60.16"""
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("D9").Formula = "=""0.372"""
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("D10").Formula = "=""56.42"""
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("D11").Formula = "=""0.0802"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("D13").Formula = "=""0.852"""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("D14").Formula = "=""22.05"""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("D20").Formula = "=""70.72"""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("D22").Formula = "=""234.95"""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("D23").Formula = "=""5.19"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("D25").Formula = "=""2.51"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("D27").Formula = "=""9.68"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("D28").Formula = "=""159.88"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("D29").Formula = "=""19.72"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("D30").Formula = "=""0.127"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("D32").Formula = "=""4.84"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("D34").Formula = "=""0.0615"""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("D36").Formula = "=""6.21"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("D39").Formula = "=""1.82"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("D40").Formula = "=""3.05"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("D42").Formula = "=""1.21"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("D46").Formula = "=""91.68"""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# --- Volume(1h) column (E) updates ---
# These values keep their surrounding whitespace, so Excel always treats
# them as plain text regardless of the numeric-looking percentage inside.
$ws.Range("E2").Value = "  -4.54%  "
$ws.Range("E3").Value = "  -4.29%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("E5").Value = "  -4.20%  "
$ws.Range("E6").Value = "  -4.24%  "
$ws.Range("E7").Value = "  -8.61%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("E10").Value = "  -5.28%  "
$ws.Range("E11").Value = "  +6.36%  "
$ws.Range("E12").Value = "  -1.49%  "
$ws.Range("E13").Value = "  -5.64%  "
$ws.Range("E14").Value = "  +7.54%  "
$ws.Range("E15").Value = "  -7.67%  "
$ws.Range("E16").Value = "  -4.17%  "
$ws.Range("E17").Value = "  -3.37%  "
$ws.Range("E18").Value = "  -4.38%  "
$ws.Range("E19").Value = "  -4.61%  "
$ws.Range("E20").Value = "  -3.50%  "
$ws.Range("E21").Value = "  -2.13%  "
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("E23").Value = "  -2.88%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E26").Value = "  -4.41%  "
$ws.Range("E27").Value = "  +1.19%  "
$ws.Range("E28").Value = "  -2.91%  "
$ws.Range("E29").Value = "  -1.03%  "
$ws.Range("E30").Value = "  +12.54%  "
$ws.Range("E31").Value = "  -2.41%  "
$ws.Range("E32").Value = "  -7.17%  "
$ws.Range("E33").Value = "  -6.03%  "
$ws.Range("E34").Value = "  +0.48%  "
$ws.Range("E35").Value = "  -7.25%  "
$ws.Range("E36").Value = "  +1.72%  "
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("E38").Value = "  -7.62%  "
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("E40").Value = "  +10.01%  "
$ws.Range("E41").Value = "  -4.13%  "
$ws.Range("E42").Value = "  -1.35%  "
$ws.Range("E43").Value = "  -3.22%  "
$ws.Range("E44").Value = "  -3.41%  "
$ws.Range("E45").Value = "  -4.97%  "
$ws.Range("E46").Value = "  -3.50%  "
$ws.Range("E47").Value = "  -5.71%  "
$ws.Range("E48").Value = "  -7.80%  "
$ws.Range("E49").Value = "  -6.41%  "
$ws.Range("E50").Value = "  -4.04%  "
$ws.Range("E51").Value = "  -4.03%  "
